$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update working hours value for B18 (row 18, Developer1 column) from 4 to 8
$ws.Range("B18").Value = 8

# Move the active cell selection to L23 (as captured in the saved view state)
$ws.Range("L23").Select()
